$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices, percentages, hour) stay as text
$ws.Range("D2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "303.18"
$ws.Range("E2").Value = "2.89%"
$ws.Range("G2").Value = "14"
$ws.Range("D3").Value = "43.04"
$ws.Range("E3").Value = "7.39%"
$ws.Range("G3").Value = "14"
$ws.Range("D4").Value = "5.052"
$ws.Range("E4").Value = "0.76%"
$ws.Range("G4").Value = "14"
$ws.Range("D5").Value = "0.07677"
$ws.Range("E5").Value = "4.84%"
$ws.Range("G5").Value = "14"
$ws.Range("D6").Value = "4.411"
$ws.Range("E6").Value = "2.25%"
$ws.Range("G6").Value = "14"
$ws.Range("D7").Value = "1.607"
$ws.Range("E7").Value = "3.75%"
$ws.Range("G7").Value = "14"
$ws.Range("D8").Value = "1.029"
$ws.Range("E8").Value = "11.30%"
$ws.Range("G8").Value = "14"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1236"
$ws.Range("E9").Value = "5.20%"
$ws.Range("G9").Value = "14"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1853"
$ws.Range("E10").Value = "2.69%"
$ws.Range("G10").Value = "14"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.09039"
$ws.Range("E11").Value = "3.36%"
$ws.Range("G11").Value = "14"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.04166"
$ws.Range("E12").Value = "-2.32%"
$ws.Range("G12").Value = "14"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.1045"
$ws.Range("E13").Value = "-0.96%"
$ws.Range("G13").Value = "14"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001269"
$ws.Range("E14").Value = "-0.51%"
$ws.Range("G14").Value = "14"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005760"
$ws.Range("E15").Value = "-1.67%"
$ws.Range("G15").Value = "14"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "0.007430"
$ws.Range("E16").Value = "1,890.09%"
$ws.Range("G16").Value = "14"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.327"
$ws.Range("E17").Value = "-0.27%"
$ws.Range("G17").Value = "14"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.384"
$ws.Range("E18").Value = "-0.24%"
$ws.Range("G18").Value = "14"
$ws.Range("D19").Value = "0.3346"
$ws.Range("E19").Value = "2.38%"
$ws.Range("G19").Value = "14"
$ws.Range("D20").Value = "8.433"
$ws.Range("E20").Value = "6.16%"
$ws.Range("G20").Value = "14"
$ws.Range("E21").Value = "1.60%"
$ws.Range("G21").Value = "14"
$ws.Range("D22").Value = "0.3192"
$ws.Range("E22").Value = "13.78%"
$ws.Range("G22").Value = "14"
$ws.Range("D23").Value = "0.04149"
$ws.Range("E23").Value = "4.61%"
$ws.Range("G23").Value = "14"
$ws.Range("D24").Value = "0.001276"
$ws.Range("E24").Value = "0.51%"
$ws.Range("G24").Value = "14"
$ws.Range("D25").Value = "0.004483"
$ws.Range("E25").Value = "17.34%"
$ws.Range("G25").Value = "14"
$ws.Range("D26").Value = "0.0001349"
$ws.Range("E26").Value = "9.27%"
$ws.Range("G26").Value = "14"
$ws.Range("G27").Value = "14"
$ws.Range("G28").Value = "14"
$ws.Range("G29").Value = "14"
$ws.Range("G30").Value = "14"
$ws.Range("G31").Value = "14"
$ws.Range("G32").Value = "14"
$ws.Range("G33").Value = "14"
$ws.Range("G34").Value = "14"
$ws.Range("G35").Value = "14"
$ws.Range("G36").Value = "14"
$ws.Range("G37").Value = "14"
$ws.Range("D38").Value = "0.02455"
$ws.Range("E38").Value = "4.04%"
$ws.Range("G38").Value = "14"
$ws.Range("D39").Value = "0.05281"
$ws.Range("E39").Value = "3.62%"
$ws.Range("G39").Value = "14"
$ws.Range("D40").Value = "0.005937"
$ws.Range("E40").Value = "-2.04%"
$ws.Range("G40").Value = "14"
$ws.Range("D41").Value = "0.007673"
$ws.Range("E41").Value = "-1.11%"
$ws.Range("G41").Value = "14"
$ws.Range("D42").Value = "0.1346"
$ws.Range("E42").Value = "3.91%"
$ws.Range("G42").Value = "14"
$ws.Range("D43").Value = "0.007358"
$ws.Range("E43").Value = "-0.34%"
$ws.Range("G43").Value = "14"
$ws.Range("D44").Value = "0.007440"
$ws.Range("E44").Value = "3.28%"
$ws.Range("G44").Value = "14"
$ws.Range("D45").Value = "0.3021"
$ws.Range("E45").Value = "2.98%"
$ws.Range("G45").Value = "14"
$ws.Range("D46").Value = "0.00006650"
$ws.Range("E46").Value = "8.46%"
$ws.Range("G46").Value = "14"
$ws.Range("E47").Value = "-0.54%"
$ws.Range("G47").Value = "14"
$ws.Range("D48").Value = "0.04478"
$ws.Range("E48").Value = "-4.30%"
$ws.Range("G48").Value = "14"
$ws.Range("E49").Value = "-0.30%"
$ws.Range("G49").Value = "14"
$ws.Range("D50").Value = "0.00002095"
$ws.Range("E50").Value = "-0.54%"
$ws.Range("G50").Value = "14"
$ws.Range("D51").Value = "0.0001996"
$ws.Range("E51").Value = "-0.54%"
$ws.Range("G51").Value = "14"
